$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = "Attr Department"
$ws.Range("I2").Value = "Attr Category"
$ws.Range("K2").Value = "Attr BrandID"

$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(11).EntireColumn.AutoFit() | Out-Null
